$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.063.56'
$ws.Range('E2').Value = '  -3.99%  '
$ws.Range('D3').Value = '3.303.19'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.64'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('B8').Value = 'LidoStakedEther'
$ws.Range('C8').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D8').Value = '3.294.58'
$ws.Range('E8').Value = '  -4.34%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.483'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.120'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('B12').Value = 'Cardano'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.407'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.882.23'
$ws.Range('E13').Value = '  -3.74%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.129'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.31'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.46%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.324.90'
$ws.Range('E16').Value = '  -3.61%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000166'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.18%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '60.110.49'
$ws.Range('E18').Value = '  -3.96%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.67'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.63%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '374.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.96%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '74.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.05%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.551'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.36%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.484.45'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000105'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.64%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.171'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.05%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.74%  '
$ws.Range('B30').Value = 'RenderToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.09%  '
$ws.Range('B31').Value = 'USDe'
$ws.Range('C31').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.62'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.81%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.41%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.29'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.77%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.18'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.51%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '166.71'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.37%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.54'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.71%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.77'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.38%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '27.54'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -13.68%  '
$ws.Range('B41').Value = 'RenzoRestakedETH'
$ws.Range('C41').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D41').Value = '3.334.34'
$ws.Range('E41').Value = '  -4.12%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0735'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.77%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.90'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.59%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.749'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.44%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.22'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.60'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.19%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.18%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '2.402.58'
$ws.Range('E48').Value = '  -6.63%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.61'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.11%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.38%  '
